# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP), G (sum) for rows 2-7.
# Column A (date) and F (Win) are unchanged by this edit.
$data = @{
    2 = @(0.127881588408715,  0.3127903958511391, 0.1575252929769615, 0.496779210170732,  1.094976487407548)
    3 = @(1.459612070389937,  10.29869402782916,  26.21740644021617,  645.3272768299601,  683.3029893683953)
    4 = @(1.459612070389937,  1.667794583268128,  26.21740644021617,  8.660232485948974,  38.00504557982321)
    5 = @(0.3048080303191223, 1.667794583268128,  337.1190423067083,  8.660232485948974,  347.7518774062445)
    6 = @(3.230985683306322,  1.667794583268128,  0.1575252929769615, 0.496779210170732,  5.553084769722144)
    7 = @(1.459612070389937,  10.29869402782916,  0.8054896365839992, 645.3272768299601,  657.8910725647631)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value2 = $vals[0]  # B: TB
    $ws.Cells.Item($row, 3).Value2 = $vals[1]  # C: d2S
    $ws.Cells.Item($row, 4).Value2 = $vals[2]  # D: K
    $ws.Cells.Item($row, 5).Value2 = $vals[3]  # E: IP
    $ws.Cells.Item($row, 7).Value2 = $vals[4]  # G: sum
}

$wb.Save()
